$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates: removed a temperature-dependent viscosity scaling term
# upstream, which shifts the downstream velocity-triangle/loss results.
$ws.Range("B2").Value = 0.13409993242852528
$ws.Range("C2").Value = 288.25533153052328
$ws.Range("D2").Value = 0.34755500670966555
$ws.Range("E2").Value = 0.99911421749194518
$ws.Range("F2").Value = -3.3585681937866774
$ws.Range("H2").Value = 716.50780749594855
$ws.Range("I2").Value = 91.550599021366125
$ws.Range("J2").Value = 92.040123781531463
$ws.Range("K2").Value = 0.59592219376540101

# Row 4 updates: loosened a3 constraint (nozzle exit swirl search bound
# E4 dropped from 105 to 93), changing the converged rotor-row solution.
$ws.Range("E4").Value = 93
$ws.Range("F4").Value = 149000.0005574617
$ws.Range("G4").Value = 0.32887959639802616
$ws.Range("H4").Value = 0.36527684904874924
$ws.Range("I4").Value = 1.3936395955527214
$ws.Range("J4").Value = 0.04852687173693504
$ws.Range("K4").Value = 0.96752073569314956

# Row 6 updates: loosened a3 constraint (bound dropped from 133 to 113),
# changing the converged final-stage solution downstream.
$ws.Range("A6").Value = -61.146699297271851
$ws.Range("B6").Value = -69.803054086655365
$ws.Range("C6").Value = -39.149212752582415
$ws.Range("D6").Value = 1.7077055259324219
$ws.Range("E6").Value = 113
$ws.Range("F6").Value = 113000.0023600851
$ws.Range("G6").Value = 0.31818390295631632
$ws.Range("H6").Value = 0.37463044360729258
$ws.Range("I6").Value = 1.7104066382790681
$ws.Range("J6").Value = 0.05035868318259603
$ws.Range("K6").Value = 0.57274607373599251
$ws.Range("L6").Value = 0.67296415604893989
$ws.Range("M6").Value = 7706451.3191588884

# Column width adjustments (column B narrows by 1 char unit; column H narrows to match column I)
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666
$ws.Columns.Item(8).ColumnWidth = 10.666666666666666

# Selection change
$ws.Range("G10").Select()
